# Update pl_mw.xlsx line-loading results for the 380 kV case
# Columns: C,D,E,F,G,H,L,M,O (idx 3,4,5,6,7,8,12,13,15); rows 2-25
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(3,4,5,6,7,8,12,13,15)

$data = @{
    2 = @(0.2143484679778282,0.02380025060995195,0.2140565575044207,0.4876855245406659,0.3305297921427695,0.4984295297481438,0.4487302040882639,10.85609145848326,1.590813011429276)
    3 = @(0.2246088004198441,0.02184503657546344,0.2004684574231206,0.4985428804552541,0.342170122640951,0.5131071354566643,0.4017053660479917,9.52247177572184,1.645694871359751)
    4 = @(0.2314175227517836,0.02064819267901186,0.1923913705296698,0.5061426335686932,0.350218186585856,0.5227868813633876,0.3730637187281616,8.700468907342611,1.682685390222474)
    5 = @(0.234317636652122,0.02016138879231733,0.1891646336448005,0.5094703732676891,0.3537196949427255,0.5268973298281594,0.3614472118009076,8.364675369643066,1.698573219334605)
    6 = @(0.2348067081310887,0.02008061100747227,0.1886326859133476,0.5100367666705807,0.3543143815443557,0.5275898264528749,0.3595215481908269,8.30886683126397,1.701260127539626)
    7 = @(0.231456129725899,0.02064162374562528,0.1923475942547341,0.5061865833515498,0.3502645169983865,0.5228416473364774,0.3729068348544899,8.695943634059404,1.68289638217837)
    8 = @(0.2177791315052531,0.02312532802681488,0.2093148146413526,0.4912330782928436,0.3343535846694863,0.5033505898925199,0.4324659760582676,10.39689593605277,1.609044327412633)
    9 = @(0.1951084218310228,0.02802568257676796,0.2448042600444467,0.4694909897461415,0.3105112337966673,0.470516029788925,0.5512590821140293,13.70905653880578,1.490984138753191)
    10 = @(0.1811412725644459,0.03164611532149308,0.2723884410755772,0.4583842861008876,0.2977739445548906,0.4498005986246056,0.6400026013581623,16.13107748827741,1.421443302518981)
    11 = @(0.1754071778944706,0.03329813417441585,0.2853027155737919,0.4544423203391617,0.2930819396432156,0.4411437846391522,0.6807532758125205,17.23117341454076,1.393735803920237)
    12 = @(0.1733280548597378,0.03392449331468583,0.2902489303474738,0.4531138299475046,0.2914691657602049,0.4379782390161751,0.6962445151754935,17.6475706522665,1.383824820013643)
    13 = @(0.1737716767427742,0.03378956033314751,0.2891811368987618,0.453392565740856,0.291809123543743,0.4386549500755876,0.6929054531267127,17.55789930200064,1.385933209525007)
    14 = @(0.1752342598171737,0.03334964927472583,0.2857085050412422,0.4543297064910732,0.2929459374037009,0.4408810829956167,0.682026519602573,17.26543402243834,1.392908676305979)
    15 = @(0.1761422439164591,0.03308029358576903,0.2835887909709953,0.4549252586784576,0.2936637893628387,0.4422593881536585,0.675370809300972,17.0862683062831,1.397257543302317)
    16 = @(0.1815287445181042,0.03153825775304853,0.2715521001509416,0.4586646367682832,0.2981032243926336,0.4503819700635887,0.637347537883187,16.05915440345979,1.423334480581957)
    17 = @(0.1849940971603843,0.03059359679125606,0.2642639736718877,0.4612462696342376,0.3011128529173561,0.4555629467638767,0.6141225264903767,15.42865915410528,1.440349299989776)
    18 = @(0.1870454036404254,0.03005072624541327,0.26010629060643,0.4628354074254801,0.3029473261187405,0.4586149339121874,0.6007996305125118,15.06585170368345,1.450504390472418)
    19 = @(0.1877498300182765,0.02986699963788197,0.258704358902861,0.4633912300414451,0.3035860435322348,0.4596605841493755,0.5962946778993796,14.94298161488325,1.454005541010844)
    20 = @(0.1846191660382601,0.030694108251339,0.2650362392729733,0.4609606281412084,0.3007817304640383,0.4550039516419275,0.6165911571426079,15.49579294282836,1.438499774222379)
    21 = @(0.1748021342613271,0.0334788403510089,0.2867269586145369,0.4540499511904699,0.2926075327756337,0.4402241388905566,0.6852202595916879,17.35134269693907,1.390843904201404)
    22 = @(0.1689254114454428,0.03530337199913447,0.3012303178564082,0.4504929578178078,0.2882236919819832,0.4312221284568096,0.7304247846223859,18.56300348761812,1.363093906954418)
    23 = @(0.1720115133061597,0.03432915106562007,0.2934585490271218,0.4523020662879347,0.2904738872137926,0.4359657271373862,0.7062643451207293,17.91639241104048,1.377588321600001)
    24 = @(0.1847884887736058,0.03064866631444829,0.2646869974477823,0.4610894401931418,0.3009311067863365,0.4552564449683985,0.6154749978049665,15.46544278085003,1.439334784121201)
    25 = @(0.2007810763769839,0.02669664225123114,0.2349503146567287,0.4745360199817696,0.3161433789587917,0.4788092065533505,0.5188840749475787,12.81526114358201,1.519967051521732)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($r, $cols[$i]).Value = $rowVals[$i]
    }
}

Write-Host "Updated pl_mw.xlsx values for case with 380 kV"